# revision_yoy_GDP.xlsx — bugfix for the naive forecaster component module.
#
# Column A previously held raw Excel date serials (e.g. 1987-12-31) formatted
# with a custom "YYYY-MM-DD HH:MM:SS" number format. The fix replaces those
# date values with plain text quarter labels ("1987Q4" ... "2024Q4") so the
# naive forecaster keys revisions by fiscal-year quarter instead of by date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (bold font, border, centered alignment,
# General number format) onto the data rows first, so the cells lose the
# custom date number format *before* we overwrite their values below. Using
# PasteSpecial (formats only) re-uses the existing header style instead of
# minting a new cell style.
$ws.Range("A1").Copy()
$ws.Range("A2:A39").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Replace each date serial with its corresponding "<year>Q4" text label.
for ($year = 1987; $year -le 2024; $year++) {
    $row = $year - 1987 + 2
    $ws.Cells.Item($row, 1).Value = "$($year)Q4"
}
